$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.947.55'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '3.074.96'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.82'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.04'
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('D9').Value = '3.068.81'
$ws.Range('E9').Value = '  -2.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.91'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000243'
$ws.Range('E13').Value = '  -2.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.55'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '3.591.90'
$ws.Range('E15').Value = '  -2.58%  '
$ws.Range('E16').Value = '  -2.17%  '
$ws.Range('D17').Value = '63.649.50'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.15'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').Value = '3.068.90'
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.06'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.59'
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.727'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.52'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('B24').Value = 'Fetch.AI'
$ws.Range('C24').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.25'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.13'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.89'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.01'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.20'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.36'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').Value = '0.0₃0852'
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  +3.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.10'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.24'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.28'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.58'
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '447.48'
$ws.Range('E42').Value = '  -1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.287'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0365'
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('D45').Value = '2.822.64'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.83'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.109'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.33'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.37'
$ws.Range('E49').Value = '  +3.58%  '
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').Value = '  -0.52%  '
